$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "249.43") are preserved verbatim as text, matching the source data.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns row by row
$ws.Cells.Item(2, 4).Value = "41.839.71"
$ws.Cells.Item(2, 5).Value = "  -1.25%  "
$ws.Cells.Item(3, 4).Value = "2.228.52"
$ws.Cells.Item(3, 5).Value = "  -0.98%  "
$ws.Cells.Item(4, 5).Value = "  -0.11%  "
$ws.Cells.Item(5, 4).Value = "249.43"
$ws.Cells.Item(5, 5).Value = "  +6.61%  "
$ws.Cells.Item(6, 5).Value = "  -0.64%  "
$ws.Cells.Item(7, 4).Value = "71.84"
$ws.Cells.Item(7, 5).Value = "  +2.89%  "
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 5).Value = "  +5.38%  "
$ws.Cells.Item(10, 4).Value = "41.39"
$ws.Cells.Item(10, 5).Value = "  +14.74%  "
$ws.Cells.Item(11, 4).Value = "0.0968"
$ws.Cells.Item(11, 5).Value = "  -2.50%  "
$ws.Cells.Item(12, 4).Value = "58.10"
$ws.Cells.Item(12, 5).Value = "  -0.54%  "
$ws.Cells.Item(13, 4).Value = "7.14"
$ws.Cells.Item(13, 5).Value = "  +5.35%  "
$ws.Cells.Item(14, 4).Value = "0.105"
$ws.Cells.Item(14, 5).Value = "  -1.29%  "
$ws.Cells.Item(15, 4).Value = "2.561.45"
$ws.Cells.Item(15, 5).Value = "  -0.85%  "
$ws.Cells.Item(16, 4).Value = "14.98"
$ws.Cells.Item(16, 5).Value = "  -1.04%  "
$ws.Cells.Item(17, 4).Value = "0.861"
$ws.Cells.Item(17, 5).Value = "  +0.11%  "
$ws.Cells.Item(18, 4).Value = "2.227.91"
$ws.Cells.Item(18, 5).Value = "  -1.00%  "
$ws.Cells.Item(19, 4).Value = "41.815.65"
$ws.Cells.Item(19, 5).Value = "  -1.10%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0966"
$ws.Cells.Item(20, 5).Value = "  -1.16%  "
$ws.Cells.Item(21, 4).Value = "6.21"
$ws.Cells.Item(21, 5).Value = "  -0.99%  "
$ws.Cells.Item(22, 4).Value = "72.92"
$ws.Cells.Item(22, 5).Value = "  -0.79%  "
$ws.Cells.Item(23, 4).Value = "235.02"
$ws.Cells.Item(23, 5).Value = "  -0.70%  "
$ws.Cells.Item(24, 4).Value = "2.12"
$ws.Cells.Item(24, 5).Value = "  +6.20%  "
$ws.Cells.Item(25, 4).Value = "4.10"
$ws.Cells.Item(25, 5).Value = "  +11.82%  "
$ws.Cells.Item(26, 4).Value = "0.999"
$ws.Cells.Item(26, 5).Value = "  -0.15%  "
$ws.Cells.Item(27, 4).Value = "2.52"
$ws.Cells.Item(27, 5).Value = "  +5.18%  "
$ws.Cells.Item(28, 4).Value = "10.65"
$ws.Cells.Item(28, 5).Value = "  +5.87%  "
$ws.Cells.Item(29, 5).Value = "  +0.06%  "
$ws.Cells.Item(30, 4).Value = "171.35"
$ws.Cells.Item(30, 5).Value = "  +1.23%  "
$ws.Cells.Item(31, 4).Value = "20.75"
$ws.Cells.Item(31, 5).Value = "  +0.59%  "
$ws.Cells.Item(32, 5).Value = "  +4.06%  "
$ws.Cells.Item(33, 5).Value = "  -1.19%  "
$ws.Cells.Item(34, 4).Value = "5.58"
$ws.Cells.Item(34, 5).Value = "  +3.43%  "
$ws.Cells.Item(35, 5).Value = "  +1.56%  "
$ws.Cells.Item(36, 4).Value = "4.73"
$ws.Cells.Item(36, 5).Value = "  +0.19%  "
$ws.Cells.Item(37, 4).Value = "26.16"
$ws.Cells.Item(37, 5).Value = "  +20.66%  "
$ws.Cells.Item(38, 4).Value = "3.95"
$ws.Cells.Item(38, 5).Value = "  +8.71%  "
$ws.Cells.Item(39, 4).Value = "0.0300"
$ws.Cells.Item(39, 5).Value = "  +11.59%  "
$ws.Cells.Item(40, 4).Value = "2.28"
$ws.Cells.Item(40, 5).Value = "  +0.57%  "
$ws.Cells.Item(41, 5).Value = "  +0.22%  "
$ws.Cells.Item(42, 4).Value = "67.31"
$ws.Cells.Item(42, 5).Value = "  +1.44%  "
$ws.Cells.Item(43, 4).Value = "12.11"
$ws.Cells.Item(43, 5).Value = "  +19.19%  "
$ws.Cells.Item(44, 4).Value = "4.99"
$ws.Cells.Item(44, 5).Value = "  +1.36%  "
$ws.Cells.Item(45, 5).Value = "  +7.02%  "
$ws.Cells.Item(46, 4).Value = "8.78"
$ws.Cells.Item(46, 5).Value = "  -2.23%  "
$ws.Cells.Item(47, 5).Value = "  -1.05%  "
$ws.Cells.Item(48, 4).Value = "4.66"
$ws.Cells.Item(48, 5).Value = "  +5.02%  "
$ws.Cells.Item(49, 5).Value = "  +0.16%  "
$ws.Cells.Item(50, 5).Value = "  +6.86%  "
$ws.Cells.Item(51, 5).Value = "  +1.07%  "

# Restore default (Normal) style on column D so no stray number-format
# style is left behind now that the values are set as text.
$dRange.Style = "Normal"

